# Add a "Ox" region-reference column (D) to the Dicionario sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D1").Value = "Ox"

# Data rows: D2:D138 -> "R1".."R137"
for ($i = 1; $i -le 137; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = "R$i"
}

# Widen the new column like the other data columns.
$ws.Columns.Item(4).ColumnWidth = 12.43

# Match the author's final selection after typing the last value.
$ws.Range("D7").Select()
